$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

$data = @(
    @("TagCollection",   "add",       1, "powerpoint-tags", "addMultipleSlideTags"),
    @("SlideCollection",  "getItemAt", 1, "powerpoint-tags", "addTagToSelectedSlide"),
    @("TagCollection",    "getItem",   1, "powerpoint-tags", "addTagToSelectedSlide"),
    @("Slide",             "delete",   1, "powerpoint-tags", "deleteSlidesByAudience"),
    @("ShapeCollection",  "getItemAt", 1, "powerpoint-tags", "addShapeTag"),
    @("TagCollection",     "delete",   1, "powerpoint-tags", "deletePresentationTag")
)

foreach ($row in $data) {
    $newRow = $lo.ListRows.Add()
    $r = $newRow.Range
    $r.Cells.Item(1, 1).Value = $row[0]
    $r.Cells.Item(1, 2).Value = $row[1]
    $r.Cells.Item(1, 3).Value = $row[2]
    $r.Cells.Item(1, 4).Value = $row[3]
    $r.Cells.Item(1, 5).Value = $row[4]
}

$ws.Range("D10").Select() | Out-Null
